$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '58.907.61'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '2.526.66'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.09'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("D9").Value = '2.526.67'
$ws.Range("E9").Value = '  +0.74%  '
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("E11").Value = '  -2.13%  '
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("D14").Value = '2.968.65'
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("D16").Value = '58.931.21'
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").Value = '2.513.00'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("E23").Value = '  +1.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.74'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.65%  '
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("E26").Value = '  -1.67%  '
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '167.55'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.72%  '
$ws.Range("E33").Value = '  +6.10%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.47'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.47%  '
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("E37").Value = '  -1.96%  '
$ws.Range("E38").Value = '  -2.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.74'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.46%  '
$ws.Range("E40").Value = '  +1.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.63'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '284.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.997'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '131.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.607'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.03%  '
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("E48").Value = '  -1.15%  '
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("E50").Value = '  -1.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.06%  '
